$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of an existing header cell (G1, "sum") onto the new
# H1 header cell so the new "Save" header matches the look of the other
# headers (bold, centered, bordered), then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill in the data values for the new "Save" column.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
